$wb = $excel.ActiveWorkbook

# --- Sheet: item_consumable ---
$wsConsumable = $wb.Worksheets.Item("item_consumable")

$wsConsumable.Range("A3").Value  = "POTION01"
$wsConsumable.Range("A4").Value  = "POTION02"
$wsConsumable.Range("A5").Value  = "POTION03"
$wsConsumable.Range("A6").Value  = "POTION04"
$wsConsumable.Range("A7").Value  = "POTION05"
$wsConsumable.Range("A8").Value  = "POTION06"
$wsConsumable.Range("A9").Value  = "POTION07"
$wsConsumable.Range("A10").Value = "POTION08"

# Narrow column A slightly to better fit the new (longer) ids.
$wsConsumable.Columns.Item(1).ColumnWidth = 7

# --- Sheet: item_equipment ---
$wsEquipment = $wb.Worksheets.Item("item_equipment")

$wsEquipment.Range("A3").Value  = "WEAPON01"
$wsEquipment.Range("A4").Value  = "WEAPON02"
$wsEquipment.Range("A5").Value  = "WEAPON03"
$wsEquipment.Range("A6").Value  = "WEAPON04"
$wsEquipment.Range("A7").Value  = "WEAPON05"
$wsEquipment.Range("A8").Value  = "WEAPON06"
$wsEquipment.Range("A9").Value  = "WEAPON07"
$wsEquipment.Range("A10").Value = "WEAPON08"

$wsEquipment.Range("A11").Value = "ARMOR01"
$wsEquipment.Range("A12").Value = "ARMOR02"
$wsEquipment.Range("A13").Value = "ARMOR03"
$wsEquipment.Range("A14").Value = "ARMOR04"
$wsEquipment.Range("A15").Value = "ARMOR05"
$wsEquipment.Range("A16").Value = "ARMOR06"
$wsEquipment.Range("A17").Value = "ARMOR07"
$wsEquipment.Range("A18").Value = "ARMOR08"

$wsEquipment.Range("A19").Value = "RING01"
$wsEquipment.Range("A20").Value = "RING02"
$wsEquipment.Range("A21").Value = "RING03"
$wsEquipment.Range("A22").Value = "RING04"
$wsEquipment.Range("A23").Value = "RING05"
$wsEquipment.Range("A24").Value = "RING06"
$wsEquipment.Range("A25").Value = "RING07"
$wsEquipment.Range("A26").Value = "RING08"
